$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Create the new "2022-Q3" sheet by duplicating "2022-Q1" (so it inherits
#    the same column layout/styles), placing it right before "2022-Q1".
# ---------------------------------------------------------------------------
$src = $wb.Worksheets.Item("2022-Q1")
$src.Copy($src)
$newSheet = $wb.Worksheets.Item("2022-Q1 (2)")
$newSheet.Name = "2022-Q3"

# Force the text-like columns to stay text (preserve leading zeros / literal
# numeric-looking strings) instead of being auto-coerced to numbers.
$newSheet.Range("B2:G3").NumberFormat = "@"

$newSheet.Range("B2").Value = "001959"
$newSheet.Range("C2").Value = "华商乐享互联灵活配置混合A"
$newSheet.Range("D2").Value = "4.62"
$newSheet.Range("E2").Value = "93.28"
$newSheet.Range("F2").Value = "2.94"
$newSheet.Range("G2").Value = "0.1358"
$newSheet.Range("H2").Value = 6

$newSheet.Range("B3").Value = "013142"
$newSheet.Range("C3").Value = "华商乐享互联灵活配置混合C"
$newSheet.Range("D3").Value = "1.08"
$newSheet.Range("E3").Value = "93.28"
$newSheet.Range("F3").Value = "2.94"
$newSheet.Range("G3").Value = "0.0318"
$newSheet.Range("H3").Value = 6

# ---------------------------------------------------------------------------
# 2) Update the "总计" (summary) sheet: add a "2022-Q3" row, keeping the
#    existing rows' data but shifted down by one. Written bottom-up so the
#    pre-existing styled cells (column A) are copied forward cleanly.
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

$total.Range("A4").Copy()
$total.Range("A5").PasteSpecial(-4122)
$total.Range("A5").Value = 3
$total.Range("B5").Value = "2020-Q4"
$total.Range("C5").Value = 3
$total.Range("D5").Value = 0.16

$total.Range("A4").Value = 2
$total.Range("B4").Value = "2021-Q4"
$total.Range("C4").Value = 3
$total.Range("D4").Value = 0.27

$total.Range("A3").Value = 1
$total.Range("B3").Value = "2022-Q1"
$total.Range("C3").Value = 2
$total.Range("D3").Value = 0.11

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q3"
$total.Range("C2").Value = 2
$total.Range("D2").Value = 0.17

# ---------------------------------------------------------------------------
# 3) Restore the originally-active tab ("2020-Q4") since editing other
#    sheets along the way shifts focus to whichever sheet was last touched.
# ---------------------------------------------------------------------------
$wb.Worksheets.Item("2020-Q4").Activate()
